# Applies the "Blind Teaser" edit:
#   - Slide 2 (Business Profile & Infrastructure): rewrite bullet text down to
#     3 bullets and remove the accompanying picture.
#   - Slide 3 (Financial & Operational Scale): rewrite bullet text down to
#     2 bullets and remove the accompanying chart.
#   - Slide 4 (Investment Highlights): rewrite bullet text down to 1 bullet
#     and remove the accompanying picture.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 - "BUSINESS PROFILE & INFRASTRUCTURE"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tb2 = $s2.Shapes.Item(5)
$tr2 = $tb2.TextFrame.TextRange

# original left/top/width/height (points) so spAutoFit recalculation doesn't
# move/resize the textbox
$L2 = $tb2.Left
$T2 = $tb2.Top
$W2 = $tb2.Width
$H2 = $tb2.Height

$tr2.Paragraphs(2).Runs(1).Text = "■ the company is an Indian electronics system design and manufacturing entity."
$tr2.Paragraphs(3).Runs(1).Text = "■ Produces subsystems, microelectronics, LCD monitors, CCTV cameras for various industries including transportation, security, automot0.9308"
# paragraphs 4,5,6 collapse down to a single paragraph (4) which then takes
# the last bullet's text; always delete the second-to-last paragraph so the
# trailing paragraph mark of the text body is never directly removed
$tr2.Paragraphs(6).Delete()
$tr2.Paragraphs(5).Delete()
$tr2.Paragraphs(4).Runs(1).Text = "■ The entity has facilities in India with a capacity of 10 million units per year."

$tb2.Left = $L2
$tb2.Top = $T2
$tb2.Width = $W2
$tb2.Height = $H2

# remove the picture next to the bullet list
$s2.Shapes.Item(6).Delete()

# ---------------------------------------------------------------------------
# Slide 3 - "FINANCIAL & OPERATIONAL SCALE"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tb3 = $s3.Shapes.Item(5)
$tr3 = $tb3.TextFrame.TextRange

$L3 = $tb3.Left
$T3 = $tb3.Top
$W3 = $tb3.Width
$H3 = $tb3.Height

$tr3.Paragraphs(2).Runs(1).Text = "■ Asset Turnover in 2023 was recorded as a value of 0.9308."
# paragraphs 3,4,5,6 collapse down to a single paragraph (3)
$tr3.Paragraphs(5).Delete()
$tr3.Paragraphs(4).Delete()
$tr3.Paragraphs(3).Delete()
$tr3.Paragraphs(3).Runs(1).Text = "■ Inventory Days, Receivable Days and Payable Days were not available for all years."

$tb3.Left = $L3
$tb3.Top = $T3
$tb3.Width = $W3
$tb3.Height = $H3

# remove the chart next to the bullet list
$s3.Shapes.Item(6).Delete()

# ---------------------------------------------------------------------------
# Slide 4 - "INVESTMENT HIGHLIGHTS"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tb4 = $s4.Shapes.Item(5)
$tr4 = $tb4.TextFrame.TextRange

$L4 = $tb4.Left
$T4 = $tb4.Top
$W4 = $tb4.Width
$H4 = $tb4.Height

# paragraphs 2,3,4,5,6 collapse down to a single paragraph (2)
$tr4.Paragraphs(5).Delete()
$tr4.Paragraphs(4).Delete()
$tr4.Paragraphs(3).Delete()
$tr4.Paragraphs(2).Delete()
$tr4.Paragraphs(2).Runs(1).Text = "■ The entity has received orders worth 50 million units in the year 2023."

$tb4.Left = $L4
$tb4.Top = $T4
$tb4.Width = $W4
$tb4.Height = $H4

# remove the picture next to the bullet list
$s4.Shapes.Item(6).Delete()
